$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 3).Value = 334
$ws.Cells.Item(2, 6).Value = 45200
$ws.Cells.Item(2, 7).Value = 30042
$ws.Cells.Item(2, 8).Value = 45231
$ws.Cells.Item(3, 5).Value = 29983
$ws.Cells.Item(3, 6).Value = 45170
$ws.Cells.Item(4, 5).Value = 29983
$ws.Cells.Item(4, 6).Value = 45170
$ws.Cells.Item(4, 7).Value = 30042
$ws.Cells.Item(4, 8).Value = 45231
$ws.Cells.Item(5, 5).Value = 30011
$ws.Cells.Item(5, 6).Value = 45200
$ws.Cells.Item(5, 7).Value = 30042
$ws.Cells.Item(5, 8).Value = 45231
$ws.Cells.Item(6, 3).Value = 442
$ws.Cells.Item(6, 6).Value = 45170
$ws.Cells.Item(6, 7).Value = 30042
$ws.Cells.Item(6, 8).Value = 45231
$ws.Cells.Item(7, 5).Value = 29952
$ws.Cells.Item(7, 6).Value = 45139
$ws.Cells.Item(7, 7).Value = 30042
$ws.Cells.Item(7, 8).Value = 45231
$ws.Cells.Item(8, 4).Value = 409
$ws.Cells.Item(8, 5).Value = 29983
$ws.Cells.Item(8, 6).Value = 45170
$ws.Cells.Item(8, 8).Value = 45231
$ws.Cells.Item(9, 5).Value = 29983
$ws.Cells.Item(9, 6).Value = 45170
$ws.Cells.Item(9, 7).Value = 30042
$ws.Cells.Item(9, 8).Value = 45231
$ws.Cells.Item(10, 3).Value = 500
$ws.Cells.Item(10, 4).Value = 482
$ws.Cells.Item(10, 5).Value = 29983
$ws.Cells.Item(10, 6).Value = 45170
$ws.Cells.Item(10, 8).Value = 45231
$ws.Cells.Item(11, 5).Value = 29983
$ws.Cells.Item(11, 6).Value = 45170
$ws.Cells.Item(11, 7).Value = 30042
$ws.Cells.Item(11, 8).Value = 45231
$ws.Cells.Item(12, 3).Value = 370
$ws.Cells.Item(12, 4).Value = 351
$ws.Cells.Item(12, 6).Value = 45170
$ws.Cells.Item(12, 8).Value = 45231
$ws.Cells.Item(13, 3).Value = 466
$ws.Cells.Item(13, 6).Value = 45170
$ws.Cells.Item(13, 7).Value = 30042
$ws.Cells.Item(13, 8).Value = 45231
$ws.Cells.Item(14, 3).Value = 423
$ws.Cells.Item(14, 4).Value = 395
$ws.Cells.Item(14, 6).Value = 45170
$ws.Cells.Item(14, 8).Value = 45231
$ws.Cells.Item(15, 3).Value = 383
$ws.Cells.Item(15, 6).Value = 45170
$ws.Cells.Item(15, 7).Value = 30011
$ws.Cells.Item(15, 8).Value = 45231
$ws.Cells.Item(16, 3).Value = 454
$ws.Cells.Item(16, 4).Value = 409
$ws.Cells.Item(16, 6).Value = 45170
$ws.Cells.Item(16, 8).Value = 45231
$ws.Cells.Item(17, 3).Value = 369
$ws.Cells.Item(17, 4).Value = 393
$ws.Cells.Item(17, 6).Value = 45170
$ws.Cells.Item(17, 8).Value = 45231
$ws.Cells.Item(18, 5).Value = 29983
$ws.Cells.Item(18, 6).Value = 45170
$ws.Cells.Item(18, 7).Value = 30042
$ws.Cells.Item(18, 8).Value = 45231
$ws.Cells.Item(19, 4).Value = 397
$ws.Cells.Item(19, 5).Value = 29983
$ws.Cells.Item(19, 6).Value = 45170
$ws.Cells.Item(19, 8).Value = 45231
$ws.Cells.Item(20, 3).Value = 480
$ws.Cells.Item(20, 6).Value = 45170
$ws.Cells.Item(20, 7).Value = 30042
$ws.Cells.Item(20, 8).Value = 45231
$ws.Cells.Item(21, 3).Value = 309
$ws.Cells.Item(21, 6).Value = 45170
$ws.Cells.Item(21, 7).Value = 30042
$ws.Cells.Item(21, 8).Value = 45231
$ws.Cells.Item(22, 3).Value = 322
$ws.Cells.Item(22, 4).Value = 366
$ws.Cells.Item(22, 6).Value = 45170
$ws.Cells.Item(22, 8).Value = 45231
$ws.Cells.Item(23, 4).Value = 314
$ws.Cells.Item(23, 5).Value = 29952
$ws.Cells.Item(23, 6).Value = 45139
$ws.Cells.Item(23, 8).Value = 45231
$ws.Cells.Item(24, 3).Value = 332
$ws.Cells.Item(24, 4).Value = 312
$ws.Cells.Item(24, 6).Value = 45170
$ws.Cells.Item(24, 8).Value = 45231
$ws.Cells.Item(25, 5).Value = 29983
$ws.Cells.Item(25, 6).Value = 45170
$ws.Cells.Item(25, 7).Value = 30042
$ws.Cells.Item(25, 8).Value = 45231
$ws.Cells.Item(26, 4).Value = 379
$ws.Cells.Item(26, 5).Value = 29983
$ws.Cells.Item(26, 6).Value = 45170
$ws.Cells.Item(26, 8).Value = 45231
$ws.Cells.Item(27, 4).Value = 216
$ws.Cells.Item(27, 5).Value = 29983
$ws.Cells.Item(27, 6).Value = 45170
$ws.Cells.Item(27, 7).Value = 35339
$ws.Cells.Item(27, 8).Value = 45231
$ws.Cells.Item(28, 3).Value = 393
$ws.Cells.Item(28, 6).Value = 45170
$ws.Cells.Item(28, 7).Value = 30042
$ws.Cells.Item(28, 8).Value = 45231
$ws.Cells.Item(29, 5).Value = 29983
$ws.Cells.Item(29, 6).Value = 45170
$ws.Cells.Item(29, 7).Value = 30042
$ws.Cells.Item(29, 8).Value = 45231
$ws.Cells.Item(30, 3).Value = 454
$ws.Cells.Item(30, 4).Value = 397
$ws.Cells.Item(30, 6).Value = 45200
$ws.Cells.Item(30, 8).Value = 45231
$ws.Cells.Item(31, 3).Value = 403
$ws.Cells.Item(31, 4).Value = 316
$ws.Cells.Item(31, 6).Value = 45139
$ws.Cells.Item(31, 8).Value = 45231
$ws.Cells.Item(32, 3).Value = 500
$ws.Cells.Item(32, 4).Value = 409
$ws.Cells.Item(32, 5).Value = 30011
$ws.Cells.Item(32, 6).Value = 45200
$ws.Cells.Item(32, 8).Value = 45231
$ws.Cells.Item(33, 3).Value = 465
$ws.Cells.Item(33, 4).Value = 316
$ws.Cells.Item(33, 6).Value = 45170
$ws.Cells.Item(33, 8).Value = 45231
$ws.Cells.Item(34, 3).Value = 358
$ws.Cells.Item(34, 4).Value = 366
$ws.Cells.Item(34, 6).Value = 45170
$ws.Cells.Item(34, 8).Value = 45231
